# 2015-2017 movie gross list: add column J ("movie page" relative URL) data + style for J2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 uses the real (non-dummy) importxml-based lookup formula, matching the source sheet
$ws.Range("J2").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("importxml(""https://www.boxofficemojo.com/yearly/chart/?view=releasedate&view2=domestic&yr=2015&sort=gross&order=DESC&p=.htm"", ""//*[@id=''body'']/table[3]/tr/td[1]/table[1]/tr/td[2]/table[1]/tr/td/table[1]/tr/td/table[1]/tr/td[2]/b/a/@href"")"),"/movies/?id=starwars7.htm")'

# J3:J101 use the generic dummy-function passthrough pattern seen throughout the sheet
$ws.Range("J3").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=jurassicpark4.htm")'
$ws.Range("J4").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=avengers2.htm")'
$ws.Range("J5").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pixar2014.htm")'
$ws.Range("J6").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=fast7.htm")'
$ws.Range("J7").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=minions.htm")'
$ws.Range("J8").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hungergames4.htm")'
$ws.Range("J9").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=scott2016.htm")'
$ws.Range("J10").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=cinderella2015.htm")'
$ws.Range("J11").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bond24.htm")'
$ws.Range("J12").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mi5.htm")'
$ws.Range("J13").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pitchperfect2.htm")'
$ws.Range("J14").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=revenant.htm")'
$ws.Range("J15").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=antman.htm")'
$ws.Range("J16").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=happysmekday.htm")'
$ws.Range("J17").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hoteltransylvania2.htm")'
$ws.Range("J18").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=fiftyshadesofgrey.htm")'
$ws.Range("J19").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=spongebob2.htm")'
$ws.Range("J20").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=straightouttacompton.htm")'
$ws.Range("J21").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=sanandreas.htm")'
$ws.Range("J22").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=furyroad.htm")'
$ws.Range("J23").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=daddyshome.htm")'
$ws.Range("J24").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=insurgent.htm")'
$ws.Range("J25").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=peanuts2015.htm")'
$ws.Range("J26").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=secretservice.htm")'
$ws.Range("J27").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pixar2013.htm")'
$ws.Range("J28").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=susancooper.htm")'
$ws.Range("J29").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=trainwreck15.htm")'
$ws.Range("J30").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=creed.htm")'
$ws.Range("J31").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=1952.htm")'
$ws.Range("J32").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=gethard.htm")'
$ws.Range("J33").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=terminator2015.htm")'
$ws.Range("J34").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=taken3.htm")'
$ws.Range("J35").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=nest.htm")'
$ws.Range("J36").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=alvin4.htm")'
$ws.Range("J37").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mazerunner2.htm")'
$ws.Range("J38").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=ted2.htm")'
$ws.Range("J39").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=goosebumps.htm")'
$ws.Range("J40").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pixels.htm")'
$ws.Range("J41").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=paddington.htm")'
$ws.Range("J42").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=intern.htm")'
$ws.Range("J43").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=coldwar2015.htm")'
$ws.Range("J44").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=paulblart2.htm")'
$ws.Range("J45").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bigshort.htm")'
$ws.Range("J46").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=warroom2015.htm")'
$ws.Range("J47").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=magicmike2.htm")'
$ws.Range("J48").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=blumhouse2015.htm")'
$ws.Range("J49").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=kevinhart15.htm")'
$ws.Range("J50").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=whiteybulger15.htm")'
$ws.Range("J51").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=newline15.htm")'
$ws.Range("J52").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=perfectguy.htm")'
$ws.Range("J53").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=joy.htm")'
$ws.Range("J54").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=fantasticfour15.htm")'
$ws.Range("J55").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=thehatefuleight.htm")'
$ws.Range("J56").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=focus2015.htm")'
$ws.Range("J57").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=southpaw2015.htm")'
$ws.Range("J58").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=insidiouschapter3.htm")'
$ws.Range("J59").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=poltergeist2015.htm")'
$ws.Range("J60").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=jupiterascending.htm")'
$ws.Range("J61").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=sicario.htm")'
$ws.Range("J62").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=uncle.htm")'
$ws.Range("J63").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=spotlight.htm")'
$ws.Range("J64").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=mcfarland.htm")'
$ws.Range("J65").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=blumhousejuly2015.htm")'
$ws.Range("J66").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=everest2015.htm")'
$ws.Range("J67").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=xmas2015.htm")'
$ws.Range("J68").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=krampus.htm")'
$ws.Range("J69").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=max2015.htm")'
$ws.Range("J70").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=ageofadaline.htm")'
$ws.Range("J71").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=brooklyn.htm")'
$ws.Range("J72").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=thelongestride.htm")'
$ws.Range("J73").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=boynextdoor.htm")'
$ws.Range("J74").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pan.htm")'
$ws.Range("J75").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=witherspoonvergara.htm")'
$ws.Range("J76").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=concussion2015.htm")'
$ws.Range("J77").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=duff.htm")'
$ws.Range("J78").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=womaningold.htm")'
$ws.Range("J79").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=bestexotic2.htm")'
$ws.Range("J80").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=cybernatural.htm")'
$ws.Range("J81").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=entourage.htm")'
$ws.Range("J82").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=papertowns.htm")'
$ws.Range("J83").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=chappie.htm")'
$ws.Range("J84").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=crimsonpeak.htm")'
$ws.Range("J85").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=walkinthewoods.htm")'
$ws.Range("J86").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=pointbreak2015.htm")'
$ws.Range("J87").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=sinister2.htm")'
$ws.Range("J88").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=lastwitchhunter.htm")'
$ws.Range("J89").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=coup.htm")'
$ws.Range("J90").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=rickiandtheflash.htm")'
$ws.Range("J91").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=womaninblack2.htm")'
$ws.Range("J92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=runallnight.htm")'
$ws.Range("J93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=lovethecoopers.htm")'
$ws.Range("J94").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=lazarus.htm")'
$ws.Range("J95").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=exmachina.htm")'
$ws.Range("J96").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=heartofthesea.htm")'
$ws.Range("J97").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=newlinehorror2015.htm")'
$ws.Range("J98").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=hitman47.htm")'
$ws.Range("J99").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=almanac.htm")'
$ws.Range("J100").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=blackorwhite.htm")'
$ws.Range("J101").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),"/movies/?id=crowe2014.htm")'

# Style J2: monospace font (Inconsolata, 11pt, black) on a solid white fill
$j2 = $ws.Range("J2")
$j2.Font.Name = "Inconsolata"
$j2.Font.Size = 11
$j2.Font.Color = 0
$j2.Interior.Pattern = 1
$j2.Interior.Color = 16777215
$j2.Interior.PatternColor = 16777215

